# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet compares two AHB format versions side by side. Its header row
# used the generic placeholder suffixes "_old" / "_new"; this rewrites them
# to the concrete format-version suffixes "_FV2310" / "_FV2404" (columns
# A-J hold the FV2310 side, column K is the unsuffixed "diff" marker,
# columns L-U hold the FV2404 side), wraps the data range in a proper
# Excel Table so the new headers double as structured-reference /
# AutoFilter column names, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells --------------------------------------------------
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Columns A-J: "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $fv2310Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

# Column K ("diff") keeps its name unchanged.

# Columns L-U: "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, 11 + 1 + $i).Value = $fv2404Headers[$i]
}

# --- 2. Wrap the data range in an Excel Table -------------------------------------
$dataRange = $ws.Range("A1:U91")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
